$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '25.377.46'
$ws.Range('E2').Value = '  -2.05%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.669.06'
$ws.Range('E3').Value = '  -4.14%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.009'
$ws.Range('E4').Value = '  +1.11%  '

# Row 5: 'BNB' -> 'BNB'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '238.73'
$ws.Range('E5').Value = '  -1.21%  '

# Row 6: 'USDC' -> 'USDC'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  +0.80%  '

# Row 7: 'XRP' -> 'XRP'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4814'
$ws.Range('E7').Value = '  -7.39%  '

# Row 8: 'Cardano' -> 'Cardano'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2624'
$ws.Range('E8').Value = '  -4.86%  '

# Row 9: 'Dogecoin' -> 'Dogecoin'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.05999'
$ws.Range('E9').Value = '  -2.69%  '

# Row 10: 'WrappedEther' -> 'TRON'
$ws.Range('B10').Value = 'TRON'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07148'
$ws.Range('E10').Value = '  -0.54%  '

# Row 11: 'TRON' -> 'WrappedEther'
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.697.82'
$ws.Range('E11').Value = '  -2.60%  '

# Row 12: 'Polygon' -> 'Polygon'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.6252'
$ws.Range('E12').Value = '  -3.11%  '

# Row 13: 'Solana' -> 'Solana'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '14.49'
$ws.Range('E13').Value = '  -3.50%  '

# Row 14: 'Polkadot' -> 'Polkadot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.631'
$ws.Range('E14').Value = '  +0.40%  '

# Row 15: 'Litecoin' -> 'Litecoin'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '73.57'
$ws.Range('E15').Value = '  -5.30%  '

# Row 16: 'Dai' -> 'BinanceUSD'
$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.010'
$ws.Range('E16').Value = '  +1.21%  '

# Row 17: 'BinanceUSD' -> 'Dai'
$ws.Range('B17').Value = 'Dai'
$ws.Range('C17').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  +0.13%  '

# Row 18: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '25.435.01'
$ws.Range('E18').Value = '  -1.83%  '

# Row 19: 'Avalanche' -> 'Avalanche'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.51'
$ws.Range('E19').Value = '  -1.83%  '

# Row 20: 'ShibaInu' -> 'ShibaInu'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000006637'
$ws.Range('E20').Value = '  -2.02%  '

# Row 21: 'Uniswap' -> 'WrappedliquidstakedEther2.0'
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.927.72'
$ws.Range('E21').Value = '  -1.92%  '

# Row 22: 'WrappedliquidstakedEther2.0' -> 'Uniswap'
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.450'
$ws.Range('E22').Value = '  +3.77%  '

# Row 23: 'Cosmos' -> 'Cosmos'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.612'
$ws.Range('E23').Value = '  -0.31%  '

# Row 24: 'Chainlink' -> 'Chainlink'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.303'
$ws.Range('E24').Value = '  +0.57%  '

# Row 25: 'Monero' -> 'Monero'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '133.73'
$ws.Range('E25').Value = '  -3.83%  '

# Row 26: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '14.84'
$ws.Range('E26').Value = '  -2.32%  '

# Row 27: 'Toncoin' -> 'Toncoin'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.387'
$ws.Range('E27').Value = '  -8.59%  '

# Row 28: 'LidoDAOToken' -> 'LidoDAOToken'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.718'
$ws.Range('E28').Value = '  -2.88%  '

# Row 29: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '102.68'
$ws.Range('E29').Value = '  -3.27%  '

# Row 30: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.838'
$ws.Range('E30').Value = '  -2.24%  '

# Row 31: 'Stellar' -> 'Stellar'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.07940'
$ws.Range('E31').Value = '  -4.25%  '

# Row 32: 'Filecoin' -> 'Filecoin'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.540'
$ws.Range('E32').Value = '  -4.57%  '

# Row 33: 'Hedera' -> 'Hedera'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04614'
$ws.Range('E33').Value = '  -0.09%  '

# Row 34: 'HuobiToken' -> 'HuobiToken'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.648'
$ws.Range('E34').Value = '  +0.20%  '

# Row 35: 'ARBITRUM' -> 'ARBITRUM'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9525'
$ws.Range('E35').Value = '  -3.83%  '

# Row 36: 'ImmutableX' -> 'ImmutableX'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.5853'
$ws.Range('E36').Value = '  -5.52%  '

# Row 37: 'MXToken' -> 'MXToken'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.654'
$ws.Range('E37').Value = '  -0.89%  '

# Row 38: 'VeChain' -> 'VeChain'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01556'
$ws.Range('E38').Value = '  -3.30%  '

# Row 39: 'TrustWalletToken' -> 'PaxDollar'
$ws.Range('B39').Value = 'PaxDollar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.006'
$ws.Range('E39').Value = '  +0.67%  '

# Row 40: 'PaxDollar' -> 'TrustWalletToken'
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.8363'
$ws.Range('E40').Value = '  +12.54%  '

# Row 41: 'RenderToken' -> 'RenderToken'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.877'
$ws.Range('E41').Value = '  -2.88%  '

# Row 42: 'Quant' -> 'Quant'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '98.93'
$ws.Range('E42').Value = '  +1.29%  '

# Row 43: 'TheSandbox' -> 'TheSandbox'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.3745'
$ws.Range('E43').Value = '  -2.89%  '

# Row 44: 'FraxShare' -> 'FraxShare'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.892'
$ws.Range('E44').Value = '  -1.87%  '

# Row 45: 'Algorand' -> 'Algorand'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.1146'
$ws.Range('E45').Value = '  +1.19%  '

# Row 46: 'Aptos' -> 'Aptos'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '6.096'
$ws.Range('E46').Value = '  -2.53%  '

# Row 47: 'Aave' -> 'Cronos'
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.05203'
$ws.Range('E47').Value = '  -0.73%  '

# Row 48: 'Cronos' -> 'Aave'
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '53.94'
$ws.Range('E48').Value = '  -1.66%  '

# Row 49: 'Elrond' -> 'Elrond'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '29.67'
$ws.Range('E49').Value = '  -2.57%  '

# Row 50: 'EnergySwap' -> 'EnergySwap'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.452'
$ws.Range('E50').Value = '  -2.12%  '

# Row 51: 'Decentraland' -> 'TrueUSD'
$ws.Range('B51').Value = 'TrueUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.007'
$ws.Range('E51').Value = '  +0.59%  '
